# Rules.xlsx update
# - Michael sheet gains two new rule rows (Foley / foleyb25@gmail.com hyperlink,
#   and Attach / attachment), matching the pattern already used on the other
#   sheets (exception handling for the "gsuite" address + an attachment rule).
# - Michael becomes the active/selected sheet (the "name pop up" tab selection
#   moves off of Julie).

$wb = $excel.ActiveWorkbook

$julie   = $wb.Worksheets.Item("Julie")
$michael = $wb.Worksheets.Item("Michael")

# Add the new rule rows to the Michael sheet.
$michael.Range("A4").Value = "Foley"
$michael.Range("B4").Value = "foleyb25@gmail.com"
$michael.Hyperlinks.Add($michael.Range("B4"), "mailto:foleyb25@gmail.com")
$michael.Range("B4").Style = "Hyperlink"

$michael.Range("A5").Value = "Attach"
$michael.Range("D5").Value = "attachment"

# Update the stored selection on Julie (no longer the active tab).
$julie.Range("D2").Select() | Out-Null

# Make Michael the active sheet/tab and set its stored selection.
$michael.Activate() | Out-Null
$michael.Range("D7").Select() | Out-Null
